$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3006201.8
$ws.Range("I98").Value = 3392.9707
$ws.Range("J98").Value = 37038030
$ws.Range("K98").Value = 3392.9707
$ws.Range("L98").Value = 37038030
$ws.Range("M98").Value = -1894.9707
$ws.Range("N98").Value = -37041026
$ws.Range("H122").Value = 3006201.8
$ws.Range("I122").Value = 3392.9707
$ws.Range("J122").Value = 37038030
$ws.Range("K122").Value = 10178.9121
$ws.Range("L122").Value = 111114090
$ws.Range("M122").Value = -7728.9121
$ws.Range("N122").Value = -111118990
$ws.Range("H132").Value = 1558
$ws.Range("I132").Value = 1334.7142
$ws.Range("J132").Value = 2674.4285
$ws.Range("K132").Value = 4004.1426
$ws.Range("L132").Value = 8023.2855
$ws.Range("M132").Value = -1474.1426
$ws.Range("N132").Value = -13083.2855
$ws.Range("H135").Value = 3792.5334
$ws.Range("I135").Value = 2966.9565
$ws.Range("J135").Value = 6505.143
$ws.Range("K135").Value = 26702.6085
$ws.Range("L135").Value = 58546.287
$ws.Range("M135").Value = -24167.6085
$ws.Range("N135").Value = -63616.287
$ws.Range("H137").Value = 1346.3077
$ws.Range("I137").Value = 921.3158
$ws.Range("J137").Value = 2499.8572
$ws.Range("K137").Value = 2763.9474
$ws.Range("L137").Value = 7499.571599999999
$ws.Range("M137").Value = -213.9474
$ws.Range("N137").Value = -12599.5716
$ws.Range("H139").Value = 78040
$ws.Range("J139").Value = 78040
$ws.Range("L139").Value = 78040
$ws.Range("N139").Value = -88320
$ws.Range("H140").Value = 97200
$ws.Range("J140").Value = 97200
$ws.Range("L140").Value = 97200
$ws.Range("N140").Value = -107560
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 643202.5600000001
$ws.Range("I32").Value = 10057.787
$ws.Range("J32").Value = 1995830
$ws.Range("K32").Value = 10057.787
$ws.Range("L32").Value = 1995830
$ws.Range("M32").Value = -9770.787
$ws.Range("N32").Value = -1996404
$ws.Range("H61").Value = 3446.4666
$ws.Range("I61").Value = 3398.7856
$ws.Range("J61").Value = 4114
$ws.Range("K61").Value = 3398.7856
$ws.Range("L61").Value = 4114
$ws.Range("M61").Value = -3186.7856
$ws.Range("N61").Value = -4538
$ws.Range("H74").Value = 1022.61536
$ws.Range("I74").Value = 1022.61536
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1022.61536
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -148.61536
$ws.Range("H77").Value = 1022.61536
$ws.Range("I77").Value = 1022.61536
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5113.0768
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -745.0767999999998
$ws.Range("H132").Value = 2654.9534
$ws.Range("I132").Value = 2289.0667
$ws.Range("J132").Value = 3499.3076
$ws.Range("K132").Value = 6867.2001
$ws.Range("L132").Value = 10497.9228
$ws.Range("M132").Value = -4337.2001
$ws.Range("N132").Value = -15557.9228
$ws.Range("H136").Value = 3446.4666
$ws.Range("I136").Value = 3398.7856
$ws.Range("J136").Value = 4114
$ws.Range("K136").Value = 10196.3568
$ws.Range("L136").Value = 12342
$ws.Range("M136").Value = -7646.356800000001
$ws.Range("N136").Value = -17442
$ws.Range("H139").Value = 38282.5
$ws.Range("J139").Value = 38282.5
$ws.Range("L139").Value = 38282.5
$ws.Range("N139").Value = -48562.5
$ws.Range("H140").Value = 103032.9
$ws.Range("J140").Value = 103032.9
$ws.Range("L140").Value = 103032.9
$ws.Range("N140").Value = -113392.9
$ws.Range("H141").Value = 65922.22
$ws.Range("J141").Value = 65922.22
$ws.Range("L141").Value = 65922.22
$ws.Range("N141").Value = -76282.22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 908.80646
$ws.Range("I94").Value = 908.80646
$ws.Range("K94").Value = 908.80646
$ws.Range("M94").Value = -457.80646
$ws.Range("H138").Value = 50750
$ws.Range("J138").Value = 50750
$ws.Range("L138").Value = 50750
$ws.Range("N138").Value = -61030
$ws.Range("H140").Value = 86726.664
$ws.Range("J140").Value = 86726.664
$ws.Range("L140").Value = 86726.664
$ws.Range("N140").Value = -97086.664
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11233.362
$ws.Range("I31").Value = 4263.346
$ws.Range("J31").Value = 15447.791
$ws.Range("K31").Value = 4263.346
$ws.Range("L31").Value = 15447.791
$ws.Range("M31").Value = -3968.346
$ws.Range("N31").Value = -16037.791
$ws.Range("H34").Value = 11233.362
$ws.Range("I34").Value = 4263.346
$ws.Range("J34").Value = 15447.791
$ws.Range("K34").Value = 4263.346
$ws.Range("L34").Value = 15447.791
$ws.Range("M34").Value = -4061.346
$ws.Range("N34").Value = -15851.791
$ws.Range("H51").Value = 9423.5
$ws.Range("J51").Value = 9423.5
$ws.Range("L51").Value = 9423.5
$ws.Range("N51").Value = -10895.5
$ws.Range("H58").Value = 1427.1305
$ws.Range("I58").Value = 841.3077
$ws.Range("J58").Value = 2188.7
$ws.Range("K58").Value = 841.3077
$ws.Range("L58").Value = 2188.7
$ws.Range("M58").Value = -638.3077
$ws.Range("N58").Value = -2594.7
$ws.Range("H61").Value = 9423.5
$ws.Range("J61").Value = 9423.5
$ws.Range("L61").Value = 9423.5
$ws.Range("N61").Value = -10119.5
$ws.Range("H68").Value = 15935.167
$ws.Range("J68").Value = 15935.167
$ws.Range("L68").Value = 15935.167
$ws.Range("N68").Value = -17433.167
$ws.Range("H71").Value = 15935.167
$ws.Range("J71").Value = 15935.167
$ws.Range("L71").Value = 47805.501
$ws.Range("N71").Value = -55293.501
$ws.Range("H122").Value = 20001528
$ws.Range("I122").Value = 1460.6666
$ws.Range("J122").Value = 28572986
$ws.Range("K122").Value = 4381.9998
$ws.Range("L122").Value = 85718958
$ws.Range("M122").Value = -1931.9998
$ws.Range("N122").Value = -85723858
$ws.Range("H136").Value = 1427.1305
$ws.Range("I136").Value = 841.3077
$ws.Range("J136").Value = 2188.7
$ws.Range("K136").Value = 2523.9231
$ws.Range("L136").Value = 6566.099999999999
$ws.Range("M136").Value = 26.07690000000002
$ws.Range("N136").Value = -11666.1
$ws.Range("H140").Value = 89900
$ws.Range("J140").Value = 89900
$ws.Range("L140").Value = 89900
$ws.Range("N140").Value = -100260
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 532.9756
$ws.Range("I122").Value = 359.96774
$ws.Range("J122").Value = 1069.3
$ws.Range("K122").Value = 3239.70966
$ws.Range("L122").Value = 9623.699999999999
$ws.Range("M122").Value = -789.70966
$ws.Range("N122").Value = -14523.7
$ws.Range("H131").Value = 9435103
$ws.Range("J131").Value = 10417656
$ws.Range("L131").Value = 31252968
$ws.Range("N131").Value = -31263048
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6472384
$ws.Range("I113").Value = 14287224
$ws.Range("J113").Value = 1001996.2
$ws.Range("K113").Value = 14287224
$ws.Range("L113").Value = 1001996.2
$ws.Range("M113").Value = -14285054
$ws.Range("N113").Value = -1006336.2
$ws.Range("H139").Value = 54075.332
$ws.Range("J139").Value = 54075.332
$ws.Range("L139").Value = 54075.332
$ws.Range("N139").Value = -64355.332
$ws.Range("H140").Value = 75756.336
$ws.Range("J140").Value = 75756.336
$ws.Range("L140").Value = 75756.336
$ws.Range("N140").Value = -86116.336
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2424.7778
$ws.Range("I61").Value = 1969.2307
$ws.Range("J61").Value = 3609.2
$ws.Range("K61").Value = 1969.2307
$ws.Range("L61").Value = 3609.2
$ws.Range("M61").Value = -1767.2307
$ws.Range("N61").Value = -4013.2
$ws.Range("H82").Value = 2560.7693
$ws.Range("I82").Value = 3687.5
$ws.Range("J82").Value = 2060
$ws.Range("K82").Value = 3687.5
$ws.Range("L82").Value = 2060
$ws.Range("M82").Value = -3326.5
$ws.Range("N82").Value = -2782
$ws.Range("H85").Value = 2560.7693
$ws.Range("I85").Value = 3687.5
$ws.Range("J85").Value = 2060
$ws.Range("K85").Value = 3687.5
$ws.Range("L85").Value = 2060
$ws.Range("M85").Value = -2439.5
$ws.Range("N85").Value = -4556
$ws.Range("H113").Value = 2424.7778
$ws.Range("I113").Value = 1969.2307
$ws.Range("J113").Value = 3609.2
$ws.Range("K113").Value = 1969.2307
$ws.Range("L113").Value = 3609.2
$ws.Range("M113").Value = 200.7692999999999
$ws.Range("N113").Value = -7949.2
$ws.Range("H138").Value = 59241.547
$ws.Range("J138").Value = 59241.547
$ws.Range("L138").Value = 59241.547
$ws.Range("N138").Value = -69521.54699999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58826764
$ws.Range("I81").Value = 2615.6924
$ws.Range("K81").Value = 5231.3848
$ws.Range("M81").Value = -4170.3848
$ws.Range("H84").Value = 58826764
$ws.Range("I84").Value = 2615.6924
$ws.Range("K84").Value = 26156.924
$ws.Range("M84").Value = -20852.924
$ws.Range("H99").Value = 27650
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 27650
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = 27650
$ws.Range("N99").Value = -33640
$ws.Range("H136").Value = 1114.7241
$ws.Range("I136").Value = 631.6087
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 1894.8261
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = 655.1739
$ws.Range("N136").Value = -14000.0001
$ws.Range("H138").Value = 87175
$ws.Range("J138").Value = 87175
$ws.Range("L138").Value = 87175
$ws.Range("N138").Value = -97455
$ws.Range("H139").Value = 57325
$ws.Range("J139").Value = 57325
$ws.Range("L139").Value = 57325
$ws.Range("N139").Value = -67605
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 79502.5
$ws.Range("J141").Value = 79502.5
$ws.Range("L141").Value = 79502.5
$ws.Range("N141").Value = -89862.5
